$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '26.505.42'
Set-TextValue 'E2' '  +0.69%  '
Set-TextValue 'D3' '1.727.67'
Set-TextValue 'E3' '  +0.56%  '
Set-TextValue 'D4' '0.9996'
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '245.24'
Set-TextValue 'E5' '  +2.44%  '
Set-TextValue 'E6' '  -0.12%  '
Set-TextValue 'D7' '0.4801'
Set-TextValue 'E7' '  +1.74%  '
Set-TextValue 'D8' '0.2662'
Set-TextValue 'E8' '  +1.35%  '
Set-TextValue 'D9' '0.06213'
Set-TextValue 'E9' '  +0.30%  '
Set-TextValue 'D10' '1.725.43'
Set-TextValue 'E10' '  +0.41%  '
Set-TextValue 'D11' '0.07149'
Set-TextValue 'E11' '  +1.24%  '
Set-TextValue 'E12' '  +1.97%  '
Set-TextValue 'D13' '0.6160'
Set-TextValue 'E13' '  +3.94%  '
Set-TextValue 'E14' '  +2.54%  '
Set-TextValue 'D15' '77.12'
Set-TextValue 'E15' '  +1.28%  '
Set-TextValue 'E16' '  -0.09%  '
Set-TextValue 'D17' '26.507.10'
Set-TextValue 'E17' '  +0.71%  '
Set-TextValue 'E18' '  -0.05%  '
Set-TextValue 'E19' '  +2.03%  '
Set-TextValue 'E20' '  +0.74%  '
Set-TextValue 'D21' '1.946.73'
Set-TextValue 'E21' '  +0.40%  '
Set-TextValue 'D22' '4.515'
Set-TextValue 'E22' '  -0.77%  '
Set-TextValue 'D23' '8.934'
Set-TextValue 'E23' '  +2.26%  '
Set-TextValue 'D24' '5.278'
Set-TextValue 'E24' '  -0.89%  '
Set-TextValue 'E25' '  +0.70%  '
Set-TextValue 'E26' '  +0.62%  '
Set-TextValue 'E27' '  +2.04%  '
Set-TextValue 'D28' '1.405'
Set-TextValue 'E28' '  -0.13%  '
Set-TextValue 'D29' '106.84'
Set-TextValue 'E29' '  -1.27%  '
Set-TextValue 'D30' '3.973'
Set-TextValue 'E30' '  -0.83%  '
Set-TextValue 'D31' '0.08021'
Set-TextValue 'E31' '  +3.81%  '
Set-TextValue 'D32' '3.706'
Set-TextValue 'E32' '  +0.59%  '
Set-TextValue 'E33' '  +2.86%  '
Set-TextValue 'E34' '  -0.10%  '
Set-TextValue 'D35' '2.616'
Set-TextValue 'E35' '  -0.01%  '
Set-TextValue 'D36' '0.6342'
Set-TextValue 'E36' '  +2.36%  '
Set-TextValue 'D37' '0.9920'
Set-TextValue 'E37' '  +1.54%  '
Set-TextValue 'D38' '0.9246'
Set-TextValue 'E38' '  -0.15%  '
Set-TextValue 'D39' '2.098'
Set-TextValue 'E39' '  +10.35%  '
Set-TextValue 'D40' '2.417'
Set-TextValue 'E40' '  +0.06%  '
Set-TextValue 'D41' '105.18'
Set-TextValue 'D42' '1.006'
Set-TextValue 'E42' '  +0.42%  '
Set-TextValue 'D43' '0.01500'
Set-TextValue 'E43' '  +1.18%  '
Set-TextValue 'D44' '5.593'
Set-TextValue 'E44' '  +4.72%  '
Set-TextValue 'D45' '0.3891'
Set-TextValue 'E45' '  +2.03%  '
Set-TextValue 'D46' '6.938'
Set-TextValue 'E46' '  +10.40%  '
Set-TextValue 'D47' '0.1182'
Set-TextValue 'E47' '  +1.69%  '
Set-TextValue 'D48' '0.05331'
Set-TextValue 'E48' '  +0.82%  '
Set-TextValue 'D49' '30.88'
Set-TextValue 'E49' '  +1.13%  '
Set-TextValue 'D50' '7.840'
Set-TextValue 'E50' '  +1.83%  '
Set-TextValue 'D51' '1.266'
Set-TextValue 'E51' '  +4.22%  '
